$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the two trailing rows (old rows 6 and 7) that were dropped from the merged dataframe
$ws.Range("A6:AY7").EntireRow.Delete()

# --- Row 2 ---
$ws.Range("A2").Value = '1975'
$ws.Range("B2").Value = 'aba roy'
$ws.Range("C2").Value = 'Marvin Barnes'
$ws.Range("D2").Value = 22
$ws.Range("E2").Value = 'SSL'
$ws.Range("G2").Value = 14
$ws.Range("H2").Value = 30
$ws.Range("I2").Value = 0.467
$ws.Range("J2").Value = 'True'
$ws.Range("K2").Value = 6345
$ws.Range("L2").Value = 1601
$ws.Range("M2").Value = '1974-75'
$ws.Range("O2").Value = 3539
$ws.Range("P2").Value = 'Marvin Barnes'
$ws.Range("Q2").Value = 1952
$ws.Range("R2").Value = 'PF'
$ws.Range("S2").Value = 23
$ws.Range("T2").Value = 1
$ws.Range("U2").Value = 'ABA'
$ws.Range("V2").Value = 'SSL'
$ws.Range("W2").Value = 77
$ws.Range("Y2").Value = 39.9
$ws.Range("Z2").Value = 10.1
$ws.Range("AA2").Value = 20.3
$ws.Range("AB2").Value = 0.498
$ws.Range("AC2").Value = 0
$ws.Range("AD2").Value = 0
$ws.Range("AE2").Value = 0
$ws.Range("AF2").Value = 10.1
$ws.Range("AG2").Value = 20.2
$ws.Range("AH2").Value = 0.499
$ws.Range("AI2").Value = 0.498
$ws.Range("AJ2").Value = 3.8
$ws.Range("AK2").Value = 5.7
$ws.Range("AL2").Value = 0.67
$ws.Range("AM2").Value = 5.4
$ws.Range("AN2").Value = 10.2
$ws.Range("AO2").Value = 15.6
$ws.Range("AP2").Value = 3.2
$ws.Range("AQ2").Value = 1.2
$ws.Range("AR2").Value = 1.8
$ws.Range("AS2").Value = 4
$ws.Range("AT2").Value = 4.3
$ws.Range("AU2").Value = 24
$ws.Range("AV2").Value = '1974-75'
$ws.Range("AW2").Value = 'Yes'
$ws.Range("AX2").Value = 'No'
$ws.Range("AY2").Value = 1975
$ws.Range("A2").Copy($ws.Range("N2"))

# --- Row 3 ---
$ws.Range("A3").Value = '1974'
$ws.Range("B3").Value = 'aba roy'
$ws.Range("C3").Value = 'Swen Nater'
$ws.Range("D3").Value = 24
$ws.Range("E3").Value = 'TOT'
$ws.Range("G3").Value = 24
$ws.Range("H3").Value = 64
$ws.Range("I3").Value = 0.375
$ws.Range("J3").Value = 'True'
$ws.Range("K3").Value = 6027
$ws.Range("L3").Value = 1541
$ws.Range("M3").Value = '1973-74'
$ws.Range("O3").Value = 4727
$ws.Range("P3").Value = 'Swen Nater'
$ws.Range("Q3").Value = 1949
$ws.Range("R3").Value = 'C'
$ws.Range("S3").Value = 25
$ws.Range("T3").Value = 1
$ws.Range("U3").Value = 'ABA'
$ws.Range("V3").Value = 'TOT'
$ws.Range("W3").Value = 79
$ws.Range("Y3").Value = 30.1
$ws.Range("Z3").Value = 5.9
$ws.Range("AA3").Value = 10.7
$ws.Range("AB3").Value = 0.552
$ws.Range("AC3").Value = 0
$ws.Range("AD3").Value = 0
$ws.Range("AE3").Value = 0
$ws.Range("AF3").Value = 5.9
$ws.Range("AG3").Value = 10.7
$ws.Range("AH3").Value = 0.553
$ws.Range("AI3").Value = 0.552
$ws.Range("AJ3").Value = 2.3
$ws.Range("AK3").Value = 3.2
$ws.Range("AL3").Value = 0.709
$ws.Range("AM3").Value = 3.6
$ws.Range("AN3").Value = 9
$ws.Range("AO3").Value = 12.6
$ws.Range("AP3").Value = 1.6
$ws.Range("AQ3").Value = 0.4
$ws.Range("AR3").Value = 0.8
$ws.Range("AS3").Value = 2.5
$ws.Range("AT3").Value = 2.7
$ws.Range("AU3").Value = 14.1
$ws.Range("AV3").Value = '1973-74'
$ws.Range("AW3").Value = 'Yes'
$ws.Range("AX3").Value = 'No'
$ws.Range("AY3").Value = 1974
$ws.Range("A3").Copy($ws.Range("N3"))

# --- Row 4 ---
$ws.Range("A4").Value = '1973'
$ws.Range("B4").Value = 'aba roy'
$ws.Range("C4").Value = 'Brian Taylor'
$ws.Range("D4").Value = 21
$ws.Range("E4").Value = 'NYA'
$ws.Range("G4").Value = 24
$ws.Range("H4").Value = 59
$ws.Range("I4").Value = 0.407
$ws.Range("J4").Value = 'True'
$ws.Range("K4").Value = 5273
$ws.Range("L4").Value = 1427
$ws.Range("M4").Value = '1972-73'
$ws.Range("O4").Value = 647
$ws.Range("P4").Value = 'Brian Taylor'
$ws.Range("Q4").Value = 1951
$ws.Range("R4").Value = 'PG'
$ws.Range("S4").Value = 22
$ws.Range("T4").Value = 1
$ws.Range("U4").Value = 'ABA'
$ws.Range("V4").Value = 'NYA'
$ws.Range("W4").Value = 63
$ws.Range("Y4").Value = 32.3
$ws.Range("Z4").Value = 6.3
$ws.Range("AA4").Value = 12.2
$ws.Range("AB4").Value = 0.515
$ws.Range("AC4").Value = 0.1
$ws.Range("AD4").Value = 0.4
$ws.Range("AE4").Value = 0.16
$ws.Range("AF4").Value = 6.2
$ws.Range("AG4").Value = 11.8
$ws.Range("AH4").Value = 0.527
$ws.Range("AI4").Value = 0.518
$ws.Range("AJ4").Value = 2.7
$ws.Range("AK4").Value = 3.6
$ws.Range("AL4").Value = 0.743
$ws.Range("AM4").Value = 1.2
$ws.Range("AN4").Value = 2
$ws.Range("AO4").Value = 3.2
$ws.Range("AP4").Value = 2.8
$ws.Range("AQ4").Value = 'N/A - Stat tracked as of the 1973-74 ABA Season'
$ws.Range("AR4").Value = 'N/A - Stat tracked as of the 1973-74 ABA Season'
$ws.Range("AS4").Value = 2.2
$ws.Range("AT4").Value = 3.5
$ws.Range("AU4").Value = 15.3
$ws.Range("AV4").Value = '1972-73'
$ws.Range("AW4").Value = 'No'
$ws.Range("AX4").Value = 'No'
$ws.Range("AY4").Value = 1973
$ws.Range("A4").Copy($ws.Range("N4"))

# --- Row 5 ---
$ws.Range("A5").Value = '1972'
$ws.Range("B5").Value = 'aba roy'
$ws.Range("C5").Value = 'Artis Gilmore'
$ws.Range("D5").Value = 22
$ws.Range("E5").Value = 'KEN'
$ws.Range("G5").Value = 38
$ws.Range("H5").Value = 60
$ws.Range("I5").Value = 0.633
$ws.Range("J5").Value = 'True'
$ws.Range("K5").Value = 4770
$ws.Range("L5").Value = 1338
$ws.Range("M5").Value = '1971-72'
$ws.Range("O5").Value = 293
$ws.Range("P5").Value = 'Artis Gilmore'
$ws.Range("Q5").Value = 1949
$ws.Range("R5").Value = 'C'
$ws.Range("S5").Value = 23
$ws.Range("T5").Value = 1
$ws.Range("U5").Value = 'ABA'
$ws.Range("V5").Value = 'KEN'
$ws.Range("W5").Value = 84
$ws.Range("Y5").Value = 43.6
$ws.Range("Z5").Value = 9.6
$ws.Range("AA5").Value = 16
$ws.Range("AB5").Value = 0.598
$ws.Range("AC5").Value = 0
$ws.Range("AD5").Value = 0
$ws.Range("AE5").Value = 0
$ws.Range("AF5").Value = 9.6
$ws.Range("AG5").Value = 16
$ws.Range("AH5").Value = 0.598
$ws.Range("AI5").Value = 0.598
$ws.Range("AJ5").Value = 4.7
$ws.Range("AK5").Value = 7.2
$ws.Range("AL5").Value = 0.646
$ws.Range("AM5").Value = 5
$ws.Range("AN5").Value = 12.7
$ws.Range("AO5").Value = 17.8
$ws.Range("AP5").Value = 2.7
$ws.Range("AQ5").Value = 'N/A - Stat tracked as of the 1973-74 ABA Season'
$ws.Range("AR5").Value = 'N/A - Stat tracked as of the 1973-74 ABA Season'
$ws.Range("AS5").Value = 4
$ws.Range("AT5").Value = 3.3
$ws.Range("AU5").Value = 23.8
$ws.Range("AV5").Value = '1971-72'
$ws.Range("AW5").Value = 'Yes'
$ws.Range("AX5").Value = 'No'
$ws.Range("AY5").Value = 1972
$ws.Range("A5").Copy($ws.Range("N5"))

